$d = $word.ActiveDocument
$d.Content.Find.Execute("sisEtema", $true, $false, $false, $false, $false,
                         $true, 1, $false, "sistema", 2)
